$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Update Riders (column C) and Average (column D) for Mon-Fri (rows 2-6)
# with the new weekly ridership run values (20161026).
$ws.Range("C2").Value = 180
$ws.Range("D2").Value = 92.49

$ws.Range("C3").Value = 171
$ws.Range("D3").Value = 93.59

$ws.Range("C4").Value = 209
$ws.Range("D4").Value = 100.35

$ws.Range("C5").Value = 168
$ws.Range("D5").Value = 97.48

$ws.Range("C6").Value = 184
$ws.Range("D6").Value = 94.38
